$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the credentials: A2 (email, hyperlinked) and B2 (password)
$ws.Range("A2").Value = "cm@ext462.com"
$ws.Range("B2").Value = "mtktB-P"

# Remove the mailto hyperlink that was attached to A2
$ws.Hyperlinks.Delete()

# Select B2 as the active cell (matches the saved selection in the sheet view)
$ws.Range("B2").Select()
